# validacion generación siniestro SISE
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the NroSiniestro value in F2 (was "1120194100370", now "1120194100385").
# Leading apostrophe keeps it entered as text (matching the cell's existing
# quote-prefixed text style) instead of Excel reinterpreting the digit string
# as a number and dropping the quote-prefix formatting.
$ws.Range("F2").Value = "'1120194100385"

# Move the active selection to F3 (as reflected in the saved view state)
$ws.Range("F3").Select()
